# Refresh the cryptos snapshot's Price (D) and Volume(1h) (E) columns
# with the latest scrape, per the GitHub Actions job that produced this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '26.275.09'
$ws.Range('E2').Value = '  -1.66%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.584.60'
$ws.Range('E3').Value = '  -1.06%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.07%  '

# Row 5: BNB
$ws.Range('E5').Value = '  -0.87%  '

# Row 6: XRP
$ws.Range('E6').Value = '  -1.44%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.06%  '

# Row 8: Dogecoin
$ws.Range('E8').Value = '  -1.08%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -0.41%  '

# Row 10: Solana
$ws.Range('D10').Value = '''19.55'
$ws.Range('E10').Value = '  -0.89%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.0845'
$ws.Range('E11').Value = '  +0.42%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.807.61'
$ws.Range('E12').Value = '  -1.05%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.586.56'
$ws.Range('E13').Value = '  -1.14%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  -0.46%  '

# Row 16: Litecoin
$ws.Range('D16').Value = '''64.38'
$ws.Range('E16').Value = '  -1.37%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '26.287.08'
$ws.Range('E17').Value = '  -1.52%  '

# Row 18: ShibaInu
$ws.Range('D18').Value = '0.0₃0740'
$ws.Range('E18').Value = '  -0.71%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''7.25'
$ws.Range('E19').Value = '  +0.37%  '

# Row 20: Dai
$ws.Range('E20').Value = '  -0.04%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''207.34'

# Row 22: Uniswap
$ws.Range('E22').Value = '  -1.14%  '

# Row 23: Toncoin
$ws.Range('E23').Value = '  -4.05%  '

# Row 24: Avalanche
$ws.Range('E24').Value = '  -1.88%  '

# Row 25: Monero
$ws.Range('E25').Value = '  +0.54%  '

# Row 26: BinanceUSD
$ws.Range('E26').Value = '  -0.05%  '

# Row 27: Cosmos
$ws.Range('D27').Value = '''7.00'
$ws.Range('E27').Value = '  -1.81%  '

# Row 28: Stellar
$ws.Range('E28').Value = '  -0.59%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '''15.29'
$ws.Range('E29').Value = '  -0.85%  '

# Row 30: Hedera
$ws.Range('E30').Value = '  -2.25%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  -0.90%  '

# Row 32: Filecoin
$ws.Range('E32').Value = '  -1.01%  '

# Row 33: WEMIXToken
$ws.Range('E33').Value = '  +13.23%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('E34').Value = '  -1.18%  '

# Row 35: Maker
$ws.Range('D35').Value = '1.283.08'
$ws.Range('E35').Value = '  -1.19%  '

# Row 36: HuobiToken
$ws.Range('E36').Value = '  +0.46%  '

# Row 37: ImmutableX
$ws.Range('D37').Value = '''0.612'
$ws.Range('E37').Value = '  +0.24%  '

# Row 38: LidoDAOToken
$ws.Range('E38').Value = '  -1.25%  '

# Row 39: VeChain
$ws.Range('E39').Value = '  -1.52%  '

# Row 40: ARBITRUM
$ws.Range('E40').Value = '  -0.66%  '

# Row 41: FraxShare
$ws.Range('E41').Value = '  +0.90%  '

# Row 43: MXToken
$ws.Range('D43').Value = '''2.13'
$ws.Range('E43').Value = '  -3.26%  '

# Row 44: Aave
$ws.Range('D44').Value = '''62.33'
$ws.Range('E44').Value = '  -1.49%  '

# Row 45: RocketPoolETH
$ws.Range('D45').Value = '1.719.76'
$ws.Range('E45').Value = '  -0.97%  '

# Row 46: Quant
$ws.Range('D46').Value = '''88.88'
$ws.Range('E46').Value = '  -2.46%  '

# Row 47: RenderToken
$ws.Range('D47').Value = '''1.55'
$ws.Range('E47').Value = '  -0.89%  '

# Row 48: Algorand
$ws.Range('E48').Value = '  +0.39%  '

# Row 49: Cronos
$ws.Range('E49').Value = '  -1.49%  '

# Row 50: USDD
$ws.Range('E50').Value = '  -0.02%  '

# Row 51: EnergySwap
$ws.Range('D51').Value = '''7.45'
$ws.Range('E51').Value = '  +0.03%  '
